# Update countries & provincias Spain
# Applies the daily COVID-19 stats refresh: updated per-country counters
# (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes) for the rows whose figures moved, including the
# three rank swaps (Libano/Mayotte, Yemen/Benin, Namibia/Jordania) where an
# updated country overtook its neighbour in the total-cases ranking, and
# bumps the "updated at" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=4;   Name='Estados Unidos';       B=3868756; C=35485; D=1787958; E=1937639; F=0; G=282; H=143159}
    @{Row=6;   Name='India';                B=1116999; C=39135; D=700367;  E=389129;  F=0; G=675; H=27503}
    @{Row=11;  Name='Chile';                B=330930;  C=2084;  D=301794;  E=20633;   F=0; G=58;  H=8503}
    @{Row=20;  Name='Alemania';             B=202822;  C=250;   D=187800;  E=5859;    F=0; G=1;   H=9163}
    @{Row=51;  Name='Barein';               B=36004;   C=0;     D=31765;   E=4113;    F=0; G=2;   H=126}
    @{Row=66;  Name='Marruecos';            B=17236;   C=221;   D=14921;   E=2042;    F=0; G=4;   H=273}
    @{Row=72;  Name='Kenia';                B=13353;   C=603;   D=5122;    E=7997;    F=0; G=9;   H=234}
    @{Row=85;  Name='Estado de Palestina';  B=8549;    C=345;   D=1921;    E=6566;    F=0; G=3;   H=62}
    @{Row=95;  Name='Mauritania';           B=5873;    C=60;    D=3436;    E=2282;    F=0; G=2;   H=155}
    @{Row=106; Name='Somalia';              B=3119;    C=8;     D=1457;    E=1569;    F=0; G=0;   H=93}
    @{Row=111; Name='Libano';               B=2859;    C=84;    D=1515;    E=1304;    F=0; G=0;   H=40}
    @{Row=112; Name='Mayotte';              B=2782;    C=0;     D=2591;    E=154;     F=0; G=0;   H=37}
    @{Row=129; Name='Sierra Leona';         B=1711;    C=10;    D=1237;    E=409;     F=0; G=0;   H=65}
    @{Row=130; Name='Yemen';                B=1606;    C=25;    D=712;     E=449;     F=0; G=2;   H=445}
    @{Row=131; Name='Benin';                B=1602;    C=0;     D=782;     E=789;     F=0; G=0;   H=31}
    @{Row=137; Name='Namibia';              B=1247;    C=44;    D=35;      E=1209;    F=0; G=1;   H=3}
    @{Row=138; Name='Jordania';             B=1218;    C=4;     D=1024;    E=183;     F=0; G=0;   H=11}
    @{Row=152; Name='Santo Tome y Principe';B=746;     C=3;     D=451;     E=281;     F=0; G=0;   H=14}
    @{Row=182; Name='Aruba';                B=113;     C=2;     D=100;     E=10;      F=0; G=0;   H=3}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}

$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 20:49"
